$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("POP")

$ws.Range("B2").Value = 1.5006
$ws.Range("C2").Value = 0.5203
$ws.Range("D2").Value = 0.0039
$ws.Range("B3").Value = -0.1614
$ws.Range("C3").Value = 0.3012
$ws.Range("D3").Value = 0.592
$ws.Range("B4").Value = 0.3008
$ws.Range("C4").Value = 0.1811
$ws.Range("D4").Value = 0.0967
$ws.Range("B5").Value = 0.14
$ws.Range("C5").Value = 0.1864
$ws.Range("D5").Value = 0.4526
$ws.Range("B6").Value = -0.3088
$ws.Range("C6").Value = 0.5905
$ws.Range("D6").Value = 0.6009
$ws.Range("B7").Value = 0.4023
$ws.Range("C7").Value = 0.1892
$ws.Range("D7").Value = 0.0335
$ws.Range("B8").Value = 0.367
$ws.Range("C8").Value = 0.1751
$ws.Range("D8").Value = 0.0361
$ws.Range("B9").Value = -0.2809
$ws.Range("C9").Value = 0.5956
$ws.Range("D9").Value = 0.6372
$ws.Range("B10").Value = 0.3463
$ws.Range("C10").Value = 0.2742
$ws.Range("D10").Value = 0.2065
$ws.Range("B11").Value = 0.1975
$ws.Range("C11").Value = 0.2795
$ws.Range("D11").Value = 0.4798
$ws.Range("B12").Value = 0.39
$ws.Range("C12").Value = 0.4221
$ws.Range("D12").Value = 0.3555
$ws.Range("B13").Value = 0.5138
$ws.Range("C13").Value = 0.4714
$ws.Range("D13").Value = 0.2758
$ws.Range("B14").Value = 0.5614
$ws.Range("C14").Value = 0.4883
$ws.Range("D14").Value = 0.2503
$ws.Range("B15").Value = 0.5573
$ws.Range("C15").Value = 0.504
$ws.Range("D15").Value = 0.2689
$ws.Range("B16").Value = 0.5876
$ws.Range("C16").Value = 0.5081
$ws.Range("D16").Value = 0.2474
$ws.Range("C17").Value = 0.5046
$ws.Range("D17").Value = 0.1967
$ws.Range("B18").Value = 0.5901
$ws.Range("C18").Value = 0.5024
$ws.Range("D18").Value = 0.2402
$ws.Range("B19").Value = 0.4991
$ws.Range("C19").Value = 0.5016
$ws.Range("D19").Value = 0.3198
$ws.Range("B20").Value = 0.5167
$ws.Range("C20").Value = 0.4917
$ws.Range("D20").Value = 0.2933
$ws.Range("B21").Value = 0.6117
$ws.Range("C21").Value = 0.4973
$ws.Range("D21").Value = 0.2186
$ws.Range("B22").Value = 0.7204
$ws.Range("C22").Value = 0.6068
$ws.Range("D22").Value = 0.2351
$ws.Range("B23").Value = -1.3215
$ws.Range("C23").Value = 0.0803
